$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), matching the style of the existing
# header row (bold font, thin border, centered/top-aligned) by copying H1's format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new I/J (I0/IF) columns for data rows 2-72.
$data = @(
    @(8,8),
    @(1,1),
    @(4,5),
    @(5,5),
    @(1,2),
    @(4,4),
    @(5,6),
    @(7,7),
    @(1,1),
    @(7,7),
    @(2,3),
    @(5,5),
    @(8,8),
    @(6,7),
    @(6,6),
    @(7,7),
    @(7,7),
    @(8,9),
    @(9,9),
    @(7,7),
    @(4,4),
    @(6,6),
    @(4,5),
    @(7,7),
    @(8,8),
    @(9,9),
    @(7,7),
    @(1,1),
    @(5,6),
    @(8,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(6,6),
    @(8,8),
    @(7,7),
    @(7,7),
    @(6,7),
    @(8,8),
    @(6,6),
    @(9,9),
    @(6,6),
    @(9,9),
    @(6,6),
    @(7,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(4,5),
    @(13,13),
    @(6,6),
    @(9,9),
    @(7,7),
    @(8,8),
    @(5,5),
    @(4,5),
    @(5,5),
    @(9,9),
    @(6,6),
    @(9,9),
    @(9,9),
    @(4,4),
    @(6,6),
    @(6,6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
